$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (Changed) date column C for rows 2 through 18:
# the serial date value 45207 (2023-10-08) becomes 45208 (2023-10-09).
$ws.Range("C2:C18").Value = 45208
